$d = $word.ActiveDocument

function Merge-IdRun([string]$idText) {
    # Build the full "<id>XXXX</id>" search text and locate it in the document.
    $full = "<id>" + $idText + "</id>"

    $rng = $d.Content
    $found = $rng.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    $start = $rng.Start
    $end = $rng.End
    $openEnd = $start + 4          # end of the literal "<id>" (4 chars)

    # Remove everything after "<id>" up to (and including) "</id>", leaving
    # only the opening "<id>" run (with its Courier New formatting) behind.
    $tail = $d.Range($openEnd, $end)
    $tail.Delete()

    # Re-insert the id value plus the closing tag right after "<id>" so the
    # new text inherits the formatting of the "<id>" run it is appended to,
    # producing a single merged run.
    $insPoint = $d.Range($openEnd, $openEnd)
    $insPoint.InsertAfter($idText + "</id>")
}

Merge-IdRun "p035v_2"
Merge-IdRun "p036r_1"
